$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing data row (648) down through the new rows (649-659)
$ws.Range("A648:V648").Copy()
$ws.Range("A649:V659").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 649
$ws.Range("A649").Value = "Entrainement"
$ws.Range("B649").Value = 45939
$ws.Range("C649").Value = "Global"
$ws.Range("D649").Value = "J-2"
$ws.Range("E649").Value = "Romain Thunet"
$ws.Range("F649").Value = "center back"
$ws.Range("G649").Value = "00:39:57"
$ws.Range("H649").Value = 2.29
$ws.Range("I649").Value = 0.22
$ws.Range("J649").Value = 2.07
$ws.Range("K649").Value = 0.13
$ws.Range("L649").Value = 0.06
$ws.Range("M649").Value = 0.03
$ws.Range("N649").Value = 0
$ws.Range("O649").Value = 2
$ws.Range("P649").Value = 3.42
$ws.Range("Q649").Value = 27.56
$ws.Range("R649").Value = 4.15
$ws.Range("S649").Value = 10
$ws.Range("T649").Value = 1
$ws.Range("U649").Value = 5
$ws.Range("V649").Value = 1

# Row 650
$ws.Range("A650").Value = "Entrainement"
$ws.Range("B650").Value = 45939
$ws.Range("C650").Value = "Global"
$ws.Range("D650").Value = "J-2"
$ws.Range("E650").Value = "Ilyes Boughanmi"
$ws.Range("F650").Value = "center forward"
$ws.Range("G650").Value = "01:10:15"
$ws.Range("H650").Value = 3.55
$ws.Range("I650").Value = 0.42
$ws.Range("J650").Value = 3.13
$ws.Range("K650").Value = 0.19
$ws.Range("L650").Value = 0.13
$ws.Range("M650").Value = 0.08
$ws.Range("N650").Value = 0.02
$ws.Range("O650").Value = 9
$ws.Range("P650").Value = 2.67
$ws.Range("Q650").Value = 32.2
$ws.Range("R650").Value = 5.81
$ws.Range("S650").Value = 24
$ws.Range("T650").Value = 15
$ws.Range("U650").Value = 15
$ws.Range("V650").Value = 13

# Row 651
$ws.Range("A651").Value = "Entrainement"
$ws.Range("B651").Value = 45939
$ws.Range("C651").Value = "Global"
$ws.Range("D651").Value = "J-2"
$ws.Range("E651").Value = "Karahali Souaré"
$ws.Range("F651").Value = "right forward"
$ws.Range("G651").Value = "01:06:43"
$ws.Range("H651").Value = 3.19
$ws.Range("I651").Value = 0.23
$ws.Range("J651").Value = 2.95
$ws.Range("K651").Value = 0.11
$ws.Range("L651").Value = 0.04
$ws.Range("M651").Value = 0.06
$ws.Range("N651").Value = 0.02
$ws.Range("O651").Value = 4
$ws.Range("P651").Value = 2.45
$ws.Range("Q651").Value = 33.5
$ws.Range("R651").Value = 4.57
$ws.Range("S651").Value = 23
$ws.Range("T651").Value = 4
$ws.Range("U651").Value = 19
$ws.Range("V651").Value = 8

# Row 652
$ws.Range("A652").Value = "Entrainement"
$ws.Range("B652").Value = 45939
$ws.Range("C652").Value = "Global"
$ws.Range("D652").Value = "J-2"
$ws.Range("E652").Value = "Amir Etien"
$ws.Range("F652").Value = "right forward"
$ws.Range("G652").Value = "01:10:34"
$ws.Range("H652").Value = 3.81
$ws.Range("I652").Value = 0.15
$ws.Range("J652").Value = 3.65
$ws.Range("K652").Value = 0.15
$ws.Range("L652").Value = 0.01
$ws.Range("M652").Value = 0
$ws.Range("N652").Value = 0
$ws.Range("O652").Value = 0
$ws.Range("P652").Value = 2.75
$ws.Range("Q652").Value = 20.58
$ws.Range("R652").Value = 4.78
$ws.Range("S652").Value = 31
$ws.Range("T652").Value = 7
$ws.Range("U652").Value = 17
$ws.Range("V652").Value = 2

# Row 653
$ws.Range("A653").Value = "Entrainement"
$ws.Range("B653").Value = 45939
$ws.Range("C653").Value = "Global"
$ws.Range("D653").Value = "J-2"
$ws.Range("E653").Value = "Kamal Bafounta"
$ws.Range("F653").Value = "center midfield"
$ws.Range("G653").Value = "01:10:54"
$ws.Range("H653").Value = 3.35
$ws.Range("I653").Value = 0.2
$ws.Range("J653").Value = 3.14
$ws.Range("K653").Value = 0.1
$ws.Range("L653").Value = 0.06
$ws.Range("M653").Value = 0.04
$ws.Range("N653").Value = 0
$ws.Range("O653").Value = 4
$ws.Range("P653").Value = 2.7
$ws.Range("Q653").Value = 27.31
$ws.Range("R653").Value = 3.81
$ws.Range("S653").Value = 11
$ws.Range("T653").Value = 0
$ws.Range("U653").Value = 4
$ws.Range("V653").Value = 0

# Row 654
$ws.Range("A654").Value = "Entrainement"
$ws.Range("B654").Value = 45939
$ws.Range("C654").Value = "Global"
$ws.Range("D654").Value = "J-2"
$ws.Range("E654").Value = "Mattheo Haon"
$ws.Range("F654").Value = "right back"
$ws.Range("G654").Value = "01:10:35"
$ws.Range("H654").Value = 3.86
$ws.Range("I654").Value = 0.28
$ws.Range("J654").Value = 3.57
$ws.Range("K654").Value = 0.1
$ws.Range("L654").Value = 0.1
$ws.Range("M654").Value = 0.07
$ws.Range("N654").Value = 0.01
$ws.Range("O654").Value = 7
$ws.Range("P654").Value = 3.19
$ws.Range("Q654").Value = 31.52
$ws.Range("R654").Value = 3.95
$ws.Range("S654").Value = 23
$ws.Range("T654").Value = 0
$ws.Range("U654").Value = 1
$ws.Range("V654").Value = 0

# Row 655
$ws.Range("A655").Value = "Entrainement"
$ws.Range("B655").Value = 45939
$ws.Range("C655").Value = "Global"
$ws.Range("D655").Value = "J-2"
$ws.Range("E655").Value = "Omar Benyounes"
$ws.Range("F655").Value = "center midfield"
$ws.Range("G655").Value = "01:09:51"
$ws.Range("H655").Value = 4.28
$ws.Range("I655").Value = 0.58
$ws.Range("J655").Value = 3.69
$ws.Range("K655").Value = 0.3
$ws.Range("L655").Value = 0.17
$ws.Range("M655").Value = 0.12
$ws.Range("N655").Value = 0
$ws.Range("O655").Value = 8
$ws.Range("P655").Value = 3.6
$ws.Range("Q655").Value = 28.52
$ws.Range("R655").Value = 4.09
$ws.Range("S655").Value = 30
$ws.Range("T655").Value = 3
$ws.Range("U655").Value = 15
$ws.Range("V655").Value = 0

# Row 656
$ws.Range("A656").Value = "Entrainement"
$ws.Range("B656").Value = 45939
$ws.Range("C656").Value = "Global"
$ws.Range("D656").Value = "J-2"
$ws.Range("E656").Value = "Yoann Martelat"
$ws.Range("F656").Value = "center midfield"
$ws.Range("G656").Value = "01:08:36"
$ws.Range("H656").Value = 3.25
$ws.Range("I656").Value = 0.22
$ws.Range("J656").Value = 3.03
$ws.Range("K656").Value = 0.1
$ws.Range("L656").Value = 0.1
$ws.Range("M656").Value = 0.03
$ws.Range("N656").Value = 0
$ws.Range("O656").Value = 3
$ws.Range("P656").Value = 2.75
$ws.Range("Q656").Value = 26.81
$ws.Range("R656").Value = 4.3
$ws.Range("S656").Value = 11
$ws.Range("T656").Value = 1
$ws.Range("U656").Value = 4
$ws.Range("V656").Value = 1

# Row 657
$ws.Range("A657").Value = "Entrainement"
$ws.Range("B657").Value = 45939
$ws.Range("C657").Value = "Global"
$ws.Range("D657").Value = "J-2"
$ws.Range("E657").Value = "Naim Ighbane"
$ws.Range("F657").Value = "center back"
$ws.Range("G657").Value = "01:09:20"
$ws.Range("H657").Value = 4.73
$ws.Range("I657").Value = 0.18
$ws.Range("J657").Value = 4.54
$ws.Range("K657").Value = 0.08
$ws.Range("L657").Value = 0.04
$ws.Range("M657").Value = 0.06
$ws.Range("N657").Value = 0
$ws.Range("O657").Value = 3
$ws.Range("P657").Value = 2.49
$ws.Range("Q657").Value = 30.07
$ws.Range("R657").Value = 4.29
$ws.Range("S657").Value = 14
$ws.Range("T657").Value = 1
$ws.Range("U657").Value = 7
$ws.Range("V657").Value = 0

# Row 658
$ws.Range("A658").Value = "Entrainement"
$ws.Range("B658").Value = 45939
$ws.Range("C658").Value = "Global"
$ws.Range("D658").Value = "J-2"
$ws.Range("E658").Value = "Jeremie Laurent"
$ws.Range("F658").Value = "left forward"
$ws.Range("G658").Value = "01:08:30"
$ws.Range("H658").Value = 3.58
$ws.Range("I658").Value = 0.45
$ws.Range("J658").Value = 3.13
$ws.Range("K658").Value = 0.2
$ws.Range("L658").Value = 0.15
$ws.Range("M658").Value = 0.09
$ws.Range("N658").Value = 0.01
$ws.Range("O658").Value = 6
$ws.Range("P658").Value = 3.03
$ws.Range("Q658").Value = 31
$ws.Range("R658").Value = 5.53
$ws.Range("S658").Value = 12
$ws.Range("T658").Value = 15
$ws.Range("U658").Value = 2
$ws.Range("V658").Value = 1

# Row 659
$ws.Range("A659").Value = "Entrainement"
$ws.Range("B659").Value = 45939
$ws.Range("C659").Value = "Global"
$ws.Range("D659").Value = "J-2"
$ws.Range("E659").Value = "Malik Boussaid"
$ws.Range("F659").Value = "right back"
$ws.Range("G659").Value = "01:10:47"
$ws.Range("H659").Value = 3.53
$ws.Range("I659").Value = 0.37
$ws.Range("J659").Value = 3.14
$ws.Range("K659").Value = 0.17
$ws.Range("L659").Value = 0.12
$ws.Range("M659").Value = 0.1
$ws.Range("N659").Value = 0
$ws.Range("O659").Value = 8
$ws.Range("P659").Value = 2.66
$ws.Range("Q659").Value = 29.08
$ws.Range("R659").Value = 5.27
$ws.Range("S659").Value = 22
$ws.Range("T659").Value = 8
$ws.Range("U659").Value = 22
$ws.Range("V659").Value = 9

# Restore the selection/active cell to match where the author left off editing
$ws.Range("D662").Select()

